$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the PrAcc document type names to their versioned (V1) labels.
$ws.Range("A15").Value = "Procurement procedure subscription V1"
$ws.Range("A16").Value = "Procurement document access V1"
$ws.Range("A17").Value = "Tender Submission V1"

# Autofit column A so its width reflects the new, longer text.
$ws.Columns.Item(1).AutoFit() | Out-Null

# Match the author's final selection state.
$ws.Range("A15:A17").Select() | Out-Null
